$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financeiro")

# Add new column header "tipo" in E1, reusing the same formatting (bold + border)
# as the existing header cells by copying D1's format onto E1.
$ws.Range("E1").Value = "tipo"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 2: Saldo inicial do dia
$ws.Range("A2").Value = "13/07/2025"
$ws.Range("C2").Value = "Saldo inicial do dia"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "entrada"

# Row 3: açogue
$ws.Range("A3").Value = "13/07/2025"
$ws.Range("C3").Value = "açogue "
$ws.Range("D3").Value = -256.73
$ws.Range("E3").Value = "saida"

# Row 4: Mercado
$ws.Range("A4").Value = "13/07/2025"
$ws.Range("C4").Value = "Mercado"
$ws.Range("D4").Value = -238.73
$ws.Range("E4").Value = "saida"
